# Add first test for rules builder class
#
# The "opponent_is_older" rule-sheet gets its first real test case: a new
# rule row ("Senior citizens don't belong in the ring") underneath the
# existing header row ("Your time is up, old timer"). The new row picks up
# the same (non-bold) formatting as the header cell above it, and both
# sheets' columns are widened slightly to fit the new content.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("rules")
$ws2 = $wb.Worksheets.Item("opponent_is_older")

# New rule/test row on the "opponent_is_older" sheet.
$ws2.Range("A2").Value = "Senior citizens don't belong in the ring"

# Match the formatting already used by A1 on that sheet (plain, no bold/underline).
$ws2.Range("A1").Copy()
$ws2.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Column auto-widening side effect of the new/longer content.
$ws2.Columns.Item(1).ColumnWidth = 31.65   # -> stored width ~32.43
$ws1.Columns.Item(2).ColumnWidth = 16.98   # -> stored width ~17.86
